$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) values are digit/period strings (e.g. '1.00', '56.921.96') that
# Excel's COM value-setter would otherwise auto-convert to numbers. Forcing the
# Text number format on each such cell immediately before assigning its Value
# keeps them stored as literal text, matching the source inline strings.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '56.921.96'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.338.14'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '528.70'
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.36'
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.337.72'
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.33'
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.349'
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.749.84'
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.41'
$ws.Range('E15').Value = '  -4.38%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '56.996.07'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.329.19'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '335.49'
$ws.Range('E19').Value = '  +1.49%  '
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.18'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '61.90'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('E26').Value = '  -3.28%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.994'
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.35'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '172.90'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0723'
$ws.Range('E31').Value = '  -3.43%  '
$ws.Range('E32').Value = '  -2.89%  '
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.929'
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.97'
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '39.18'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.77'
$ws.Range('E41').Value = '  +8.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '148.28'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.375'
$ws.Range('E43').Value = '  -3.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.59'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '282.98'
$ws.Range('E45').Value = '  -2.91%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0931'
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.74'
$ws.Range('E48').Value = '  +3.01%  '
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('E50').Value = '  +6.56%  '
$ws.Range('E51').Value = '  -1.51%  '
